$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells, styled like the existing bold/centered headers ---
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G1").Value = "MSE_median"
$ws.Range("H1").Value = "MAE_median"
$ws.Range("I1").Value = "Dir_accuracy"

# --- Row 2 (RandomForestRegressor): new metric values ---
$ws.Range("G2").Value = 0.000421530831398454
$ws.Range("H2").Value = 0.02053118935964799
$ws.Range("I2").Value = 0.451063829787234

# --- Row 3 (Naive): new metric values; Dir_accuracy (I3) stays blank, like F3 ---
$ws.Range("G3").Value = 0.001083194070471167
$ws.Range("H3").Value = 0.03291191380748265

$ws.Range("F3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (materializes the blank cell)
$excel.CutCopyMode = 0
